$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.443.49'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '2.451.14'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '583.32'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('D6').Value = '144.51'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').Value = '2.448.55'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '0.108'
$ws.Range('E10').Value = '  -3.55%  '
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '5.22'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('E13').Value = '  -3.33%  '
$ws.Range('D14').Value = '26.62'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').Value = '0.0000174'
$ws.Range('E15').Value = '  -3.49%  '
$ws.Range('D16').Value = '2.905.35'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '62.282.06'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').Value = '2.451.26'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '10.93'
$ws.Range('E19').Value = '  -3.25%  '
$ws.Range('D20').Value = '7.16'
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('D21').Value = '330.78'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('E23').Value = '  -3.37%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '66.12'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('D26').Value = '9.49'
$ws.Range('E26').Value = '  +5.69%  '
$ws.Range('D27').Value = '630.87'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('D28').Value = '0.0₃0964'
$ws.Range('E28').Value = '  -5.88%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -4.31%  '
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('D35').Value = '4.95'
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('E37').Value = '  -6.03%  '
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('D40').Value = '149.90'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('D41').Value = '18.39'
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D43').Value = '42.48'
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '2.49'
$ws.Range('E45').Value = '  -4.25%  '
$ws.Range('D46').Value = '144.08'
$ws.Range('E46').Value = '  -3.17%  '
$ws.Range('E47').Value = '  -3.11%  '
$ws.Range('D48').Value = '0.0527'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('D49').Value = '0.603'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('D50').Value = '19.70'
$ws.Range('E50').Value = '  -7.29%  '
$ws.Range('E51').Value = '  +8.13%  '
